$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 data ---
$ws.Range("C2").Value = "Pineda"
$ws.Range("D2").Value = "Juan"
$ws.Range("G2").Value = 1023935749
$ws.Range("L2").Value = "calle 19 1 10"

# --- Update row 3 data ---
$ws.Range("B3").Value = "Vergara"
$ws.Range("C3").Value = "Hernandez"
$ws.Range("D3").Value = "Jeimy"
$ws.Range("J3").Value = "prueba@uniandes.edu.co"
$ws.Range("L3").Value = "calle 19 1 10"

# --- Turn J3 into a mailto hyperlink (adds Hyperlink style + relationship) ---
$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:prueba@uniandes.edu.co")

# --- Update the saved selection to match the new active cell ---
[void]$ws.Range("F2").Select()
